$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" rows (16-18) get reshuffled: the period that was
# "1806" (row 16) becomes "1804", and the period that was "1804" (row 18)
# becomes "1806". The middle row (1805) is unchanged. Each period keeps
# its own "Valor Mora" value, so the F column values effectively swap
# between row 16 and row 18 as well.
$ws.Range("E16").Value = "1804"
$ws.Range("F16").Value = 31249

$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 19791
